$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Todo ")
$ws.Activate()

# New todo items: "cd" and "gcd", both with PIC/Creator "fish"
$ws.Range("A30").Value = "cd"
$ws.Range("A31").Value = "gcd"

$ws.Range("D30").Value = "fish"
$ws.Range("E30").Value = "fish"
$ws.Range("D31").Value = "fish"
$ws.Range("E31").Value = "fish"

$ws.Range("E31").Select()
